# Scheduled-runner data refresh for the per-job market/profit tables.
# Columns H:N on each class sheet hold market pricing + computed leve
# profit figures (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ],
# LeveProfit[NQ/HQ]); this run repulls them from the market data source.
# A few rows also gain or lose their trailing LeveProfitHQ (N) cell when
# an HQ price becomes available/unavailable for that snapshot.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2922.923
$ws.Range("J51").Value = 3299.6
$ws.Range("L51").Value = 3299.6
$ws.Range("N51").Value = -4267.6
$ws.Range("H64").Value = 5336.1035
$ws.Range("I64").Value = 3476.2104
$ws.Range("J64").Value = 8869.9
$ws.Range("K64").Value = 3476.2104
$ws.Range("L64").Value = 8869.9
$ws.Range("M64").Value = -3228.2104
$ws.Range("N64").Value = -9365.9
$ws.Range("H67").Value = 5336.1035
$ws.Range("I67").Value = 3476.2104
$ws.Range("J67").Value = 8869.9
$ws.Range("K67").Value = 3476.2104
$ws.Range("L67").Value = 8869.9
$ws.Range("M67").Value = -2618.2104
$ws.Range("N67").Value = -10585.9
$ws.Range("H99").Value = 599.75
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("H129").Value = 76925870
$ws.Range("I129").Value = 90911256
$ws.Range("J129").Value = 6263
$ws.Range("K129").Value = 272733768
$ws.Range("L129").Value = 18789
$ws.Range("M129").Value = -272728768
$ws.Range("N129").Value = -28789
$ws.Range("H132").Value = 2058.9048
$ws.Range("I132").Value = 1992.9744
$ws.Range("K132").Value = 5978.9232
$ws.Range("M132").Value = -3448.9232
$ws.Range("H134").Value = 95076.664
$ws.Range("J134").Value = 95076.664
$ws.Range("L134").Value = 95076.664
$ws.Range("N134").Value = -105216.664
$ws.Range("H138").Value = 6947356
$ws.Range("J138").Value = 9094268
$ws.Range("L138").Value = 27282804
$ws.Range("N138").Value = -27293084
$ws.Range("H141").Value = 2282.6
$ws.Range("I141").Value = 2282.6
$ws.Range("K141").Value = 6847.799999999999
$ws.Range("M141").Value = -1667.799999999999
$ws.Range("N99").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 931.5294
$ws.Range("I2").Value = 694.8461
$ws.Range("K2").Value = 694.8461
$ws.Range("M2").Value = -581.8461
$ws.Range("H32").Value = 8117.9326
$ws.Range("I32").Value = 3602.459
$ws.Range("J32").Value = 29305.924
$ws.Range("K32").Value = 3602.459
$ws.Range("L32").Value = 29305.924
$ws.Range("M32").Value = -3315.459
$ws.Range("N32").Value = -29879.924
$ws.Range("H61").Value = 3705.6792
$ws.Range("I61").Value = 1840.2593
$ws.Range("J61").Value = 5642.846
$ws.Range("K61").Value = 1840.2593
$ws.Range("L61").Value = 5642.846
$ws.Range("M61").Value = -1628.2593
$ws.Range("N61").Value = -6066.846
$ws.Range("H116").Value = 931.5294
$ws.Range("I116").Value = 694.8461
$ws.Range("K116").Value = 694.8461
$ws.Range("M116").Value = 1599.1539
$ws.Range("H122").Value = 1948.6666
$ws.Range("I122").Value = 973
$ws.Range("J122").Value = 3900
$ws.Range("K122").Value = 2919
$ws.Range("L122").Value = 11700
$ws.Range("M122").Value = -469
$ws.Range("N122").Value = -16600
$ws.Range("H132").Value = 3547.861
$ws.Range("I132").Value = 3341.3447
$ws.Range("J132").Value = 4403.4287
$ws.Range("K132").Value = 10024.0341
$ws.Range("L132").Value = 13210.2861
$ws.Range("M132").Value = -7494.034100000001
$ws.Range("N132").Value = -18270.2861
$ws.Range("H136").Value = 3705.6792
$ws.Range("I136").Value = 1840.2593
$ws.Range("J136").Value = 5642.846
$ws.Range("K136").Value = 5520.7779
$ws.Range("L136").Value = 16928.538
$ws.Range("M136").Value = -2970.7779
$ws.Range("N136").Value = -22028.538

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 931.5294
$ws.Range("I3").Value = 694.8461
$ws.Range("K3").Value = 694.8461
$ws.Range("M3").Value = -580.8461
$ws.Range("H134").Value = 1534.25
$ws.Range("I134").Value = 1523.0294
$ws.Range("J134").Value = 1725
$ws.Range("K134").Value = 4569.0882
$ws.Range("L134").Value = 5175
$ws.Range("M134").Value = -2034.0882
$ws.Range("N134").Value = -10245

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35639.71
$ws.Range("I31").Value = 50682.57
$ws.Range("J31").Value = 4049.7
$ws.Range("K31").Value = 50682.57
$ws.Range("L31").Value = 4049.7
$ws.Range("M31").Value = -50387.57
$ws.Range("N31").Value = -4639.7
$ws.Range("H34").Value = 35639.71
$ws.Range("I34").Value = 50682.57
$ws.Range("J34").Value = 4049.7
$ws.Range("K34").Value = 50682.57
$ws.Range("L34").Value = 4049.7
$ws.Range("M34").Value = -50480.57
$ws.Range("N34").Value = -4453.7
$ws.Range("H99").Value = 2785.1428
$ws.Range("I99").Value = 2785.1428
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2785.1428
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1287.1428
$ws.Range("H126").Value = 2785.1428
$ws.Range("I126").Value = 2785.1428
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8355.428400000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5885.428400000001
$ws.Range("H130").Value = 35000
$ws.Range("J130").Value = 35000
$ws.Range("L130").Value = 35000
$ws.Range("N130").Value = -45040
$ws.Range("H132").Value = 2073.1765
$ws.Range("I132").Value = 2073.1765
$ws.Range("K132").Value = 6219.529500000001
$ws.Range("M132").Value = -3689.529500000001
$ws.Range("H134").Value = 9812.807000000001
$ws.Range("I134").Value = 5786.625
$ws.Range("J134").Value = 23616.857
$ws.Range("K134").Value = 17359.875
$ws.Range("L134").Value = 70850.571
$ws.Range("M134").Value = -14824.875
$ws.Range("N134").Value = -75920.571
$ws.Range("N99").ClearContents()
$ws.Range("N126").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 233.64
$ws.Range("I12").Value = 152
$ws.Range("J12").Value = 265.3889
$ws.Range("K12").Value = 456
$ws.Range("L12").Value = 796.1667
$ws.Range("M12").Value = -283
$ws.Range("N12").Value = -1142.1667
$ws.Range("H68").Value = 5555944
$ws.Range("I68").Value = 415.83334
$ws.Range("K68").Value = 1247.50002
$ws.Range("M68").Value = -436.5000199999999
$ws.Range("H71").Value = 5555944
$ws.Range("I71").Value = 415.83334
$ws.Range("K71").Value = 3742.50006
$ws.Range("M71").Value = 313.4999399999997
$ws.Range("H80").Value = 1568.125
$ws.Range("J80").Value = 1790.8334
$ws.Range("L80").Value = 5372.5002
$ws.Range("N80").Value = -7244.5002
$ws.Range("H83").Value = 1568.125
$ws.Range("J83").Value = 1790.8334
$ws.Range("L83").Value = 16117.5006
$ws.Range("N83").Value = -25477.5006
$ws.Range("H141").Value = 114449.664
$ws.Range("I141").Value = 1002.3333
$ws.Range("K141").Value = 3006.9999
$ws.Range("M141").Value = 2173.0001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 3255.2
$ws.Range("I41").Value = 3499.5
$ws.Range("J41").Value = 3092.3333
$ws.Range("K41").Value = 3499.5
$ws.Range("L41").Value = 3092.3333
$ws.Range("M41").Value = -3144.5
$ws.Range("N41").Value = -3802.3333
$ws.Range("H102").Value = 90910110
$ws.Range("I102").Value = 1052
$ws.Range("K102").Value = 1052
$ws.Range("M102").Value = 570
$ws.Range("H132").Value = 6260.5386
$ws.Range("I132").Value = 2892.6667
$ws.Range("K132").Value = 8678.000100000001
$ws.Range("M132").Value = -6148.000100000001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2643.25
$ws.Range("I22").Value = 1924.75
$ws.Range("J22").Value = 2882.75
$ws.Range("K22").Value = 1924.75
$ws.Range("L22").Value = 2882.75
$ws.Range("M22").Value = -1629.75
$ws.Range("N22").Value = -3472.75
$ws.Range("H27").Value = 2643.25
$ws.Range("I27").Value = 1924.75
$ws.Range("J27").Value = 2882.75
$ws.Range("K27").Value = 1924.75
$ws.Range("L27").Value = 2882.75
$ws.Range("M27").Value = -1817.75
$ws.Range("N27").Value = -3096.75
$ws.Range("H31").Value = 2603.7058
$ws.Range("I31").Value = 186.5
$ws.Range("J31").Value = 6056.857
$ws.Range("K31").Value = 186.5
$ws.Range("L31").Value = 6056.857
$ws.Range("M31").Value = 61.5
$ws.Range("N31").Value = -6552.857
$ws.Range("H46").Value = 1369.1428
$ws.Range("I46").Value = 1199.6
$ws.Range("J46").Value = 1793
$ws.Range("K46").Value = 1199.6
$ws.Range("L46").Value = 1793
$ws.Range("M46").Value = -1011.6
$ws.Range("N46").Value = -2169
$ws.Range("H55").Value = 163.45454
$ws.Range("I55").Value = 174.52942
$ws.Range("K55").Value = 174.52942
$ws.Range("M55").Value = -1.529419999999988
$ws.Range("H122").Value = 4250.769
$ws.Range("J122").Value = 4390.154
$ws.Range("L122").Value = 13170.462
$ws.Range("N122").Value = -18070.462
$ws.Range("H132").Value = 3401.9143
$ws.Range("I132").Value = 3147.4
$ws.Range("J132").Value = 4929
$ws.Range("K132").Value = 9442.200000000001
$ws.Range("L132").Value = 14787
$ws.Range("M132").Value = -6912.200000000001
$ws.Range("N132").Value = -19847

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 223190.6
$ws.Range("J34").Value = 28988.25
$ws.Range("L34").Value = 28988.25
$ws.Range("N34").Value = -29394.25
$ws.Range("H38").Value = 4283.3335
$ws.Range("I38").Value = 2925
$ws.Range("J38").Value = 7000
$ws.Range("K38").Value = 2925
$ws.Range("L38").Value = 7000
$ws.Range("M38").Value = -2452
$ws.Range("N38").Value = -7946
$ws.Range("H39").Value = 21000
$ws.Range("J39").Value = 21000
$ws.Range("L39").Value = 21000
$ws.Range("N39").Value = -21826
$ws.Range("H40").Value = 22200
$ws.Range("J40").Value = 22200
$ws.Range("L40").Value = 22200
$ws.Range("N40").Value = -22498
$ws.Range("H42").Value = 22025
$ws.Range("J42").Value = 22025
$ws.Range("L42").Value = 22025
$ws.Range("N42").Value = -22781
$ws.Range("H113").Value = 523.71875
$ws.Range("I113").Value = 534.64
$ws.Range("K113").Value = 1603.92
$ws.Range("M113").Value = 566.0799999999999
$ws.Range("H131").Value = 147499.67
$ws.Range("J131").Value = 147499.67
$ws.Range("L131").Value = 147499.67
$ws.Range("N131").Value = -157579.67
$ws.Range("H132").Value = 4006.6924
$ws.Range("J132").Value = 3980.4
$ws.Range("L132").Value = 11941.2
$ws.Range("N132").Value = -17001.2
